# EPBDS-9677 Support dependency on rules from tests in OpenL Maven Plugin
#
# Rewrites the Project1-Main "Test sayHello" smoke-test sheet into a
# Spreadsheet-style rule "spr" (with a Steps/Formula step table) while
# keeping the Environment/dependency block further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous layout had two merged "title" cells (C3:D3 and C13:D13).
# Drop the merges and wipe the whole previous block clean before writing
# the new content so no stale cells/styles linger.
$ws.Range("C3:D3").UnMerge()
$ws.Range("C13:D13").UnMerge()
$ws.Range("C3:D14").Clear()

# New spreadsheet-rule signature.
$ws.Range("C4").Value = "Spreadsheet SpreadsheetResult spr(String name, Integer age)"

# Steps / Formula table.
$ws.Range("C5").Value = "Steps"
$ws.Range("D5").Value = "Formula"

$ws.Range("C6").Value = "Step1"
$ws.Range("C7").Value = "Step2"

# These two hold formula-looking text that must stay literal text (quote
# prefix), not be evaluated as real formulas.
$ws.Range("D6").Value = "'= sayHello(name)"
$ws.Range("D7").Value = "'= ""I am "" + age + "" age old."""

# Environment / dependency block (unchanged content, shifted down one row).
$ws.Range("C14").Value = "Environment"
$ws.Range("C15").Value = "dependency"
$ws.Range("D15").Value = "Project2-*"
